$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RUDE (B2): was numeric 123412342134 -> now text "123412342134 " (trailing space).
# Force text type so the trailing space and leading digits aren't re-parsed as a number.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "123412342134 "

# EXP (D2): new value "OR"
$ws.Range("D2").Value = "OR"

# SEXO (J2): "M" -> "H"
$ws.Range("J2").Value = "H"

# DEP. NAC (N2): "NINGUNO" -> "TARIJA"
$ws.Range("N2").Value = "TARIJA"

# PROV. NAC. (O2): "NINGUNO" -> "ARCE"
$ws.Range("O2").Value = "ARCE"
